$wb = $excel.ActiveWorkbook

# Add the new weekly sheet after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = 'magapoke_2025-11-05'

# Header row (bold, centered, bordered -- matches the other weekly sheets)
$newSheet.Range('A1').Value = 'rank'
$newSheet.Range('B1').Value = 'title'
$headerRange = $newSheet.Range("A1:B1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Match the page margins used by the other weekly sheets
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Ranking rows
$newSheet.Cells.Item(2, 1).Value = 1
$newSheet.Cells.Item(2, 2).Value = '信じていた仲間達にダンジョン奥地で殺されかけたがギフト『無限ガチャ』でレベル9999の仲間達を手に入れて元パーティーメンバーと世界に復讐＆『ざまぁ！』します！'
$newSheet.Cells.Item(3, 1).Value = 2
$newSheet.Cells.Item(3, 2).Value = 'ギルティサークル'
$newSheet.Cells.Item(4, 1).Value = 3
$newSheet.Cells.Item(4, 2).Value = '宇宙兄弟'
$newSheet.Cells.Item(5, 1).Value = 4
$newSheet.Cells.Item(5, 2).Value = 'みいちゃんと山田さん'
$newSheet.Cells.Item(6, 1).Value = 5
$newSheet.Cells.Item(6, 2).Value = 'ドラハチ'
$newSheet.Cells.Item(7, 1).Value = 6
$newSheet.Cells.Item(7, 2).Value = '薫る花は凛と咲く'
$newSheet.Cells.Item(8, 1).Value = 7
$newSheet.Cells.Item(8, 2).Value = '島耕作'
$newSheet.Cells.Item(9, 1).Value = 8
$newSheet.Cells.Item(9, 2).Value = 'イレギュラーズ'
$newSheet.Cells.Item(10, 1).Value = 9
$newSheet.Cells.Item(10, 2).Value = '転生したら第七王子だったので、気ままに魔術を極めます'
$newSheet.Cells.Item(11, 1).Value = 10
$newSheet.Cells.Item(11, 2).Value = '愛妻の裏アカ'
$newSheet.Cells.Item(12, 1).Value = 11
$newSheet.Cells.Item(12, 2).Value = '君が僕らを悪魔と呼んだ頃'
$newSheet.Cells.Item(13, 1).Value = 12
$newSheet.Cells.Item(13, 2).Value = '魔術ギルド総帥～生まれ変わって今更やり直す2度目の学院生活～'
$newSheet.Cells.Item(14, 1).Value = 13
$newSheet.Cells.Item(14, 2).Value = '十字架のろくにん'
$newSheet.Cells.Item(15, 1).Value = 14
$newSheet.Cells.Item(15, 2).Value = '黄昏町プリズナーズ'
$newSheet.Cells.Item(16, 1).Value = 15
$newSheet.Cells.Item(16, 2).Value = '黒猫と魔女の教室'
$newSheet.Cells.Item(17, 1).Value = 16
$newSheet.Cells.Item(17, 2).Value = 'ハードワーカー中田'
$newSheet.Cells.Item(18, 1).Value = 17
$newSheet.Cells.Item(18, 2).Value = 'WIND BREAKER'
$newSheet.Cells.Item(19, 1).Value = 18
$newSheet.Cells.Item(19, 2).Value = '魔女と傭兵'
$newSheet.Cells.Item(20, 1).Value = 19
$newSheet.Cells.Item(20, 2).Value = 'ひゃくえむ。'
$newSheet.Cells.Item(21, 1).Value = 20
$newSheet.Cells.Item(21, 2).Value = '味方が弱すぎて補助魔法に徹していた宮廷魔法師、追放されて最強を目指す'
$newSheet.Cells.Item(22, 1).Value = 21
$newSheet.Cells.Item(22, 2).Value = '転生貴族、鑑定スキルで成り上がる～弱小領地を受け継いだので、優秀な人材を増やしていたら、最強領地になってた～'
$newSheet.Cells.Item(23, 1).Value = 22
$newSheet.Cells.Item(23, 2).Value = 'GALAXIAS'
$newSheet.Cells.Item(24, 1).Value = 23
$newSheet.Cells.Item(24, 2).Value = 'K-9~警視庁公安部公安第9課異能対策係~'
$newSheet.Cells.Item(25, 1).Value = 24
$newSheet.Cells.Item(25, 2).Value = 'となりの黒川さん'
$newSheet.Cells.Item(26, 1).Value = 25
$newSheet.Cells.Item(26, 2).Value = 'グラぱらっ！'
$newSheet.Cells.Item(27, 1).Value = 26
$newSheet.Cells.Item(27, 2).Value = '蒼く染めろ'
$newSheet.Cells.Item(28, 1).Value = 27
$newSheet.Cells.Item(28, 2).Value = 'アルキメデスの大戦'
$newSheet.Cells.Item(29, 1).Value = 28
$newSheet.Cells.Item(29, 2).Value = '異世界ウォーキング'
$newSheet.Cells.Item(30, 1).Value = 29
$newSheet.Cells.Item(30, 2).Value = 'せいぶつ部の田辺くん'
$newSheet.Cells.Item(31, 1).Value = 30
$newSheet.Cells.Item(31, 2).Value = '南海トラフ巨大地震'
$newSheet.Cells.Item(32, 1).Value = 31
$newSheet.Cells.Item(32, 2).Value = '幼馴染とはラブコメにならない'
$newSheet.Cells.Item(33, 1).Value = 32
$newSheet.Cells.Item(33, 2).Value = '降り積もれ孤独な死よ'
$newSheet.Cells.Item(34, 1).Value = 33
$newSheet.Cells.Item(34, 2).Value = 'ナキナギ'
$newSheet.Cells.Item(35, 1).Value = 34
$newSheet.Cells.Item(35, 2).Value = 'デッドアカウント'
$newSheet.Cells.Item(36, 1).Value = 35
$newSheet.Cells.Item(36, 2).Value = 'ハナバス　苔石花江のバスケ論'
$newSheet.Cells.Item(37, 1).Value = 36
$newSheet.Cells.Item(37, 2).Value = 'Aランクパーティを離脱した俺は、元教え子たちと迷宮深部を目指す。'
$newSheet.Cells.Item(38, 1).Value = 37
$newSheet.Cells.Item(38, 2).Value = 'さわらないで小手指くん'
$newSheet.Cells.Item(39, 1).Value = 38
$newSheet.Cells.Item(39, 2).Value = '食糧人類-Starving Anonymous-'
$newSheet.Cells.Item(40, 1).Value = 39
$newSheet.Cells.Item(40, 2).Value = 'FAIRY TAIL 100 YEARS QUEST'
$newSheet.Cells.Item(41, 1).Value = 40
$newSheet.Cells.Item(41, 2).Value = '追放された転生王子、『自動製作《オートクラフト》』スキルで領地を爆速で開拓し最強の村を作ってしまう〜最強クラフトスキルで始める、楽々領地開拓スローライフ〜'
$newSheet.Cells.Item(42, 1).Value = 41
$newSheet.Cells.Item(42, 2).Value = '阿武ノーマル'
$newSheet.Cells.Item(43, 1).Value = 42
$newSheet.Cells.Item(43, 2).Value = '皇女転生　～伝説の大魔導士（♂）、姫騎士となりて伝説の令嬢騎士団を作り無双する～'
$newSheet.Cells.Item(44, 1).Value = 43
$newSheet.Cells.Item(44, 2).Value = '時々ボソッとロシア語でデレる隣のアーリャさん'
$newSheet.Cells.Item(45, 1).Value = 44
$newSheet.Cells.Item(45, 2).Value = '触手魔術師の成り上がり'
$newSheet.Cells.Item(46, 1).Value = 45
$newSheet.Cells.Item(46, 2).Value = '不遇職【鑑定士】が実は最強だった～奈落で鍛えた最強の【神眼】で無双する～'
$newSheet.Cells.Item(47, 1).Value = 46
$newSheet.Cells.Item(47, 2).Value = '可愛いだけじゃない式守さん'
$newSheet.Cells.Item(48, 1).Value = 47
$newSheet.Cells.Item(48, 2).Value = 'ジュミドロ'
$newSheet.Cells.Item(49, 1).Value = 48
$newSheet.Cells.Item(49, 2).Value = 'Destiny Unchain Online 〜吸血鬼少女となって、やがて『赤の魔王』と呼ばれるようになりました〜'
$newSheet.Cells.Item(50, 1).Value = 49
$newSheet.Cells.Item(50, 2).Value = '屋根の下のアルテミス'
$newSheet.Cells.Item(51, 1).Value = 50
$newSheet.Cells.Item(51, 2).Value = 'おやすみ ふみさん'
$newSheet.Cells.Item(52, 1).Value = 51
$newSheet.Cells.Item(52, 2).Value = 'アオバノバスケ'
$newSheet.Cells.Item(53, 1).Value = 52
$newSheet.Cells.Item(53, 2).Value = '限界集落を脱村した錬金術士、都会で"最強"なのがバレまくる。～老害どもにはいい加減愛想が尽きました～'
$newSheet.Cells.Item(54, 1).Value = 53
$newSheet.Cells.Item(54, 2).Value = 'ルックスＹを選んでしまいました 〜やり込んでいるゲームに転生したはずなのに、未実装のガチャで攻略をすることになった件〜'
$newSheet.Cells.Item(55, 1).Value = 54
$newSheet.Cells.Item(55, 2).Value = 'ストーカー行為がバレて人生終了男'
$newSheet.Cells.Item(56, 1).Value = 55
$newSheet.Cells.Item(56, 2).Value = '冰剣の魔術師が世界を統べる〜世界最強の魔術師である少年は、魔術学院に入学する〜'
$newSheet.Cells.Item(57, 1).Value = 56
$newSheet.Cells.Item(57, 2).Value = '恋ニ非ズ'
$newSheet.Cells.Item(58, 1).Value = 57
$newSheet.Cells.Item(58, 2).Value = '我間乱 ―修羅―'
$newSheet.Cells.Item(59, 1).Value = 58
$newSheet.Cells.Item(59, 2).Value = '春くらり'
$newSheet.Cells.Item(60, 1).Value = 59
$newSheet.Cells.Item(60, 2).Value = 'いじめるヤバイ奴'
$newSheet.Cells.Item(61, 1).Value = 60
$newSheet.Cells.Item(61, 2).Value = '辺境の薬師、都でSランク冒険者となる～英雄村の少年がチート薬で無自覚無双〜'
$newSheet.Cells.Item(62, 1).Value = 61
$newSheet.Cells.Item(62, 2).Value = 'この世界がいずれ滅ぶことを、俺だけが知っている～モンスターが現れた世界で、死に戻りレベルアップ～'
$newSheet.Cells.Item(63, 1).Value = 62
$newSheet.Cells.Item(63, 2).Value = '東京卍リベンジャーズ～場地圭介からの手紙～'
$newSheet.Cells.Item(64, 1).Value = 63
$newSheet.Cells.Item(64, 2).Value = '念願の悪役令嬢（ラスボス）の身体を手に入れたぞ！'
$newSheet.Cells.Item(65, 1).Value = 64
$newSheet.Cells.Item(65, 2).Value = 'お母さん冒険者、ログインボーナスでスキル【主婦】に目覚めました。週一貰えるチラシで冒険者生活頑張ります！'
$newSheet.Cells.Item(66, 1).Value = 65
$newSheet.Cells.Item(66, 2).Value = 'ブルーロック'
$newSheet.Cells.Item(67, 1).Value = 66
$newSheet.Cells.Item(67, 2).Value = 'ヒロインは絶望しました。'
$newSheet.Cells.Item(68, 1).Value = 67
$newSheet.Cells.Item(68, 2).Value = '東京ネオンスキャンダル'
$newSheet.Cells.Item(69, 1).Value = 68
$newSheet.Cells.Item(69, 2).Value = '勇者と呼ばれた後に　―そして無双男は家族を創る―'
$newSheet.Cells.Item(70, 1).Value = 69
$newSheet.Cells.Item(70, 2).Value = '不遇職『鍛冶師』だけど最強です ～気づけば何でも作れるようになっていた男ののんびりスローライフ～'
$newSheet.Cells.Item(71, 1).Value = 70
$newSheet.Cells.Item(71, 2).Value = '放課後、ぼくは君になる'
$newSheet.Cells.Item(72, 1).Value = 71
$newSheet.Cells.Item(72, 2).Value = '四刀流の最強配信者～やり込んだVRゲームの設定が現実世界に反映されたので、廃止予定だった戦闘職で無双します～'
$newSheet.Cells.Item(73, 1).Value = 72
$newSheet.Cells.Item(73, 2).Value = '母という呪縛 娘という牢獄'
$newSheet.Cells.Item(74, 1).Value = 73
$newSheet.Cells.Item(74, 2).Value = 'リスナーに騙されてダンジョンの最下層から脱出RTAすることになった'
$newSheet.Cells.Item(75, 1).Value = 74
$newSheet.Cells.Item(75, 2).Value = 'シャングリラ・フロンティア～クソゲーハンター、神ゲーに挑まんとす～'
$newSheet.Cells.Item(76, 1).Value = 75
$newSheet.Cells.Item(76, 2).Value = '異世界グルメで成り上がり無双～山に追放されたので、のんびりキャンプを楽しんでいたらいつの間にか強くなっていて、王侯貴族や実力者たちが俺を放っておいてくれません。一方、俺を追放した貴族たちは破滅が始まる～'
$newSheet.Cells.Item(77, 1).Value = 76
$newSheet.Cells.Item(77, 2).Value = 'なれの果ての僕ら'
$newSheet.Cells.Item(78, 1).Value = 77
$newSheet.Cells.Item(78, 2).Value = '魁の花巫女'
$newSheet.Cells.Item(79, 1).Value = 78
$newSheet.Cells.Item(79, 2).Value = '剣帝学院の魔眼賢者'
$newSheet.Cells.Item(80, 1).Value = 79
$newSheet.Cells.Item(80, 2).Value = 'ともだちづくり'
$newSheet.Cells.Item(81, 1).Value = 80
$newSheet.Cells.Item(81, 2).Value = 'イジらないで、長瀞さん'
$newSheet.Cells.Item(82, 1).Value = 81
$newSheet.Cells.Item(82, 2).Value = 'デスティニーラバーズ'
$newSheet.Cells.Item(83, 1).Value = 82
$newSheet.Cells.Item(83, 2).Value = 'お願い、脱がシて。'
$newSheet.Cells.Item(84, 1).Value = 83
$newSheet.Cells.Item(84, 2).Value = 'ダメスキル【自動機能】が覚醒しました～あれ、ギルドのスカウトの皆さん、俺を「いらない」って言ってませんでした？～'
$newSheet.Cells.Item(85, 1).Value = 84
$newSheet.Cells.Item(85, 2).Value = '劣等人の魔剣使い　スキルボードを駆使して最強に至る'
$newSheet.Cells.Item(86, 1).Value = 85
$newSheet.Cells.Item(86, 2).Value = '五輪の女神さま 〜なでしこ寮のメダルごはん〜'
$newSheet.Cells.Item(87, 1).Value = 86
$newSheet.Cells.Item(87, 2).Value = 'DAYS外伝'
$newSheet.Cells.Item(88, 1).Value = 87
$newSheet.Cells.Item(88, 2).Value = 'お嬢様の僕'
$newSheet.Cells.Item(89, 1).Value = 88
$newSheet.Cells.Item(89, 2).Value = 'それがメイドのカンナです'
$newSheet.Cells.Item(90, 1).Value = 89
$newSheet.Cells.Item(90, 2).Value = 'はっちぽっちぱんち'
$newSheet.Cells.Item(91, 1).Value = 90
$newSheet.Cells.Item(91, 2).Value = '最弱な僕は＜壁抜けバグ＞で成り上がる～壁をすり抜けたら、初回クリア報酬を無限回収できました！～'
$newSheet.Cells.Item(92, 1).Value = 91
$newSheet.Cells.Item(92, 2).Value = '卒業アルバムの彼女たち'
$newSheet.Cells.Item(93, 1).Value = 92
$newSheet.Cells.Item(93, 2).Value = '金田一少年の事件簿外伝 犯人たちの事件簿'
$newSheet.Cells.Item(94, 1).Value = 93
$newSheet.Cells.Item(94, 2).Value = 'MYS'
$newSheet.Cells.Item(95, 1).Value = 94
$newSheet.Cells.Item(95, 2).Value = '「無能はいらない」と言われたから絶縁してやった　～最強の四天王に育てられた俺は、冒険者となり無双する～'
$newSheet.Cells.Item(96, 1).Value = 95
$newSheet.Cells.Item(96, 2).Value = 'はじめの一歩'
$newSheet.Cells.Item(97, 1).Value = 96
$newSheet.Cells.Item(97, 2).Value = '英雄と魔女の転生ラブコメ'
$newSheet.Cells.Item(98, 1).Value = 97
$newSheet.Cells.Item(98, 2).Value = '鉱石令嬢〜没落した悪役令嬢が炭鉱で一山当てるまでのお話〜'
$newSheet.Cells.Item(99, 1).Value = 98
$newSheet.Cells.Item(99, 2).Value = '人間消失'
$newSheet.Cells.Item(100, 1).Value = 99
$newSheet.Cells.Item(100, 2).Value = '追放されなかった男　～二度目の人生は土下座から始まりました～'
$newSheet.Cells.Item(101, 1).Value = 100
$newSheet.Cells.Item(101, 2).Value = '彼女、お借りします'
